# Swap the data of row 2 and row 3 (all columns A:AY) to match the target workbook state.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric columns ---
$ws.Range("A2").Value = 106077283
$ws.Range("A3").Value = 104986863
$ws.Range("B2").Value = 56278
$ws.Range("B3").Value = 57193
$ws.Range("E2").Value = 100011
$ws.Range("E3").Value = 206004
$ws.Range("Q2").Value = 537773.3909779217
$ws.Range("Q3").Value = 537888.8853063835
$ws.Range("R2").Value = 6668679.681769322
$ws.Range("R3").Value = 6669232.05540918
$ws.Range("S2").Value = 10
$ws.Range("S3").Value = 50

# --- Text columns that are safe from auto-type-conversion (plain strings, times) ---
$ws.Range("F2").Value = "Kungsörn"
$ws.Range("F3").Value = "Skogshare"
$ws.Range("G2").Value = "Aquila chrysaetos"
$ws.Range("G3").Value = "Lepus timidus"
$ws.Range("H2").Value = "(Linnaeus, 1758)"
$ws.Range("H3").Value = "Linnaeus, 1758"
$ws.Range("M2").Value = "förbiflygande"
$ws.Range("M3").Value = "gående/springande"
$ws.Range("P2").Value = "Styggtjärnsberget, Dlr"
$ws.Range("P3").Value = "Gläfse, Jörken, Dlr"
$ws.Range("Z2").Value = "00:00"
$ws.Range("Z3").Value = "13:30"
$ws.Range("AB2").Value = "00:00"
$ws.Range("AB3").Value = "13:30"
$ws.Range("AW2").Value = "Samuel Keith"
$ws.Range("AW3").Value = "Lars Mattsson"
$ws.Range("AX2").Value = "Samuel Keith"
$ws.Range("AX3").Value = "Lars Mattsson"

# --- Date-like text columns: force Text format first to avoid auto date conversion ---
$ws.Range("Y2:Y3").NumberFormat = "@"
$ws.Range("AA2:AA3").NumberFormat = "@"
$ws.Range("Y2").Value = "2022-03-09"
$ws.Range("Y3").Value = "2022-12-07"
$ws.Range("AA2").Value = "2022-03-09"
$ws.Range("AA3").Value = "2022-12-07"

# --- Columns where a cell needs to be created on one row and cleared on the other ---
# (K, L, N move from row 2 to row 3; AC moves from row 3 to row 2)
$ws.Range("K2").ClearContents()
$ws.Range("K3").NumberFormat = "@"
$ws.Range("K3").Value = ""
$ws.Range("L2").ClearContents()
$ws.Range("L3").NumberFormat = "@"
$ws.Range("L3").Value = ""
$ws.Range("N2").ClearContents()
$ws.Range("N3").NumberFormat = "@"
$ws.Range("N3").Value = ""
$ws.Range("AC2").Value = "Mötte örn 1 och de kollade in varandra och tog några svängar tillsammans i två omgångar. Uppfattades av mig som uppvaktning, men ej att de var ett etablerat par. De skiljdes sedan och denna örn drog vidare åt väster. Åldern bedömdes t subad me"
$ws.Range("AC3").ClearContents()
